$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 9, 13, 19, 25, 31, 37, 42, 48, 52, 58, 64, 70, 76)

foreach ($r in $rows) {
    foreach ($col in @("I", "J", "K", "L", "M", "N")) {
        $ws.Range("$col$r").Value = 1
    }
    foreach ($col in @("O", "P", "Q", "R", "S", "T")) {
        $ws.Range("$col$r").Value = 0
    }
}
